# The commit adds an ExcelDataReader test case: the amount in D1 is now
# supplied as a formatted string ("1,20,000") instead of a plain number,
# and the view is scrolled/selected so D1 is in focus.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 becomes a text value -> new shared string, and the SUM(B1:D1) formula
# in F1 recalculates (text is ignored by SUM, so it drops from 13097 to 710).
$ws.Range("D1").Value = "1,20,000"

# Reflect the updated selection/scroll position from the saved view state.
$ws.Range("D1").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

$wb.Save()
